$d = $word.ActiveDocument

# The "Termina prestiti" bullet currently reads:
#   "Termina prestiti: selezionando questa voce l'utente puo annullare
#    tutti i suoi prestiti, facendo tornare le risorse associate in
#    archivio."
# and has the (hidden) "_GoBack" bookmark sitting between "Termina" and
# " prestiti". The edit rewords the sentence (adding the option to
# cancel a single loan instead of all of them) and the bookmark ends up
# relocated to the very end of the paragraph, after the final period.

# 1. Remove the existing "_GoBack" bookmark; it will be re-created at
#    the end of the paragraph once the text has been rewritten.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Insert "scegliere se " right before "annullare tutti i suoi prestiti".
$rng = $d.Content
$rng.Find.Execute("annullare tutti i suoi prestiti", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.Start, $rng.Start)
$insertPoint.InsertBefore("scegliere se ")

# 3. Insert " oppure selezionarne solamente uno" right after
#    "annullare tutti i suoi prestiti" (before the following comma).
$rng = $d.Content
$rng.Find.Execute("annullare tutti i suoi prestiti", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.End, $rng.End)
$insertPoint.InsertAfter(" oppure selezionarne solamente uno")

# 4. Insert " ai prestiti" right after "facendo tornare le risorse associate"
#    (before " in archivio").
$rng = $d.Content
$rng.Find.Execute("facendo tornare le risorse associate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.End, $rng.End)
$insertPoint.InsertAfter(" ai prestiti")

# 5. Re-add the "_GoBack" bookmark at the very end of the paragraph,
#    right after the trailing period but before the paragraph mark.
#    Adding a zero-length bookmark exactly at the paragraph-mark offset
#    is unreliable in this engine, so a temporary marker character is
#    placed right before the paragraph mark, the bookmark is inserted
#    just ahead of that marker, and the marker is then removed again.
$rng = $d.Content
$rng.Find.Execute("Termina prestiti", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)
$paraEnd = $para.Range.End - 1
$marker = $d.Range($paraEnd, $paraEnd)
$marker.InsertBefore("@@MARK@@")

$rng2 = $d.Content
$rng2.Find.Execute("@@MARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPoint = $d.Range($rng2.Start, $rng2.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint)

$rng3 = $d.Content
$rng3.Find.Execute("@@MARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Delete()
